# Updated status; changes to intron/exon
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (INTRON/EXON task) - status + notes updated
$ws.Range("B8").Value = "Does probably not need a function"
$ws.Range("C8").Value = "Does probably not need a function."
$ws.Range("D8").Value = "Work in progress"
$ws.Range("E8").Value = "Work in progress"

# The much shorter notes text no longer needs the tall row; match the
# author's re-wrapped row height.
$ws.Rows.Item(8).RowHeight = 30

# Row 15 (RESTRICTION ENZYMES task) - argument/notes updated
$ws.Range("B15").Value = "DNA sequence and exons positions"
$ws.Range("C15").Value = "1 string and 1 hash"
$ws.Range("D15").Value = "Task different from what we thought. Done. Please se module ezymes.pm"

# Update the view selection to match the author's saved cursor position
$ws.Range("D11").Select()
